$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A2:A4").Value = "Complete"

$ws.Range("F12").Select()
